$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns for rows 2-51.
# Price values are forced to text (NumberFormat "@") then restored to their
# original (default) style, since several of them parse as plain numbers
# (e.g. "1.009") and Excel would otherwise auto-convert them to numeric values.

$dCell = $ws.Range("D2")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "28.530.85"
$dCell.Style = $dStyle
$ws.Range("E2").Value = "  +0.46%  "

$dCell = $ws.Range("D3")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.872.62"
$dCell.Style = $dStyle
$ws.Range("E3").Value = "  -0.26%  "

$dCell = $ws.Range("D4")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.009"
$dCell.Style = $dStyle
$ws.Range("E4").Value = "  -0.60%  "

$dCell = $ws.Range("D5")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "315.64"
$dCell.Style = $dStyle
$ws.Range("E5").Value = "  -0.04%  "

$dCell = $ws.Range("D6")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.008"
$dCell.Style = $dStyle
$ws.Range("E6").Value = "  -0.53%  "

$dCell = $ws.Range("D7")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.5069"
$dCell.Style = $dStyle
$ws.Range("E7").Value = "  -1.26%  "

$dCell = $ws.Range("D8")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.3892"
$dCell.Style = $dStyle
$ws.Range("E8").Value = "  -0.76%  "

$dCell = $ws.Range("D9")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.08353"
$dCell.Style = $dStyle
$ws.Range("E9").Value = "  +0.56%  "

$dCell = $ws.Range("D10")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "41.76"
$dCell.Style = $dStyle
$ws.Range("E10").Value = "  -0.39%  "

$dCell = $ws.Range("D11")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.101"
$dCell.Style = $dStyle
$ws.Range("E11").Value = "  -1.71%  "

$dCell = $ws.Range("D12")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "6.207"
$dCell.Style = $dStyle
$ws.Range("E12").Value = "  -0.99%  "

$dCell = $ws.Range("D13")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.872.32"
$dCell.Style = $dStyle
$ws.Range("E13").Value = "  -2.62%  "

$dCell = $ws.Range("D14")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "20.38"
$dCell.Style = $dStyle
$ws.Range("E14").Value = "  +0.06%  "

$dCell = $ws.Range("D15")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "7.229"
$dCell.Style = $dStyle
$ws.Range("E15").Value = "  -0.33%  "

$dCell = $ws.Range("D16")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.010"
$dCell.Style = $dStyle
$ws.Range("E16").Value = "  -0.51%  "

$dCell = $ws.Range("D17")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001102"
$dCell.Style = $dStyle
$ws.Range("E17").Value = "  -0.35%  "

$dCell = $ws.Range("D18")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "91.06"
$dCell.Style = $dStyle
$ws.Range("E18").Value = "  -0.31%  "

$dCell = $ws.Range("D19")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.06703"
$dCell.Style = $dStyle
$ws.Range("E19").Value = "  -0.32%  "

$dCell = $ws.Range("D20")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "17.67"
$dCell.Style = $dStyle
$ws.Range("E20").Value = "  -0.35%  "

$dCell = $ws.Range("D21")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.007"
$dCell.Style = $dStyle
$ws.Range("E21").Value = "  -0.60%  "

$dCell = $ws.Range("D22")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "5.920"
$dCell.Style = $dStyle
$ws.Range("E22").Value = "  -1.22%  "

$dCell = $ws.Range("D23")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "28.529.11"
$dCell.Style = $dStyle
$ws.Range("E23").Value = "  +0.31%  "

$dCell = $ws.Range("D24")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "11.06"
$dCell.Style = $dStyle
$ws.Range("E24").Value = "  -0.89%  "

$dCell = $ws.Range("D25")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.232"
$dCell.Style = $dStyle
$ws.Range("E25").Value = "  -1.33%  "

$dCell = $ws.Range("D26")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.086.38"
$dCell.Style = $dStyle
$ws.Range("E26").Value = "  -2.35%  "

$dCell = $ws.Range("D27")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "161.89"
$dCell.Style = $dStyle
$ws.Range("E27").Value = "  +0.61%  "

$dCell = $ws.Range("D28")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "20.61"
$dCell.Style = $dStyle
$ws.Range("E28").Value = "  -0.92%  "

$dCell = $ws.Range("D29")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.346"
$dCell.Style = $dStyle
$ws.Range("E29").Value = "  -3.58%  "

$dCell = $ws.Range("D30")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "125.99"
$dCell.Style = $dStyle
$ws.Range("E30").Value = "  -0.21%  "

$dCell = $ws.Range("D31")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.1043"
$dCell.Style = $dStyle
$ws.Range("E31").Value = "  -2.16%  "

$dCell = $ws.Range("D32")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.036"
$dCell.Style = $dStyle
$ws.Range("E32").Value = "  -0.91%  "

$dCell = $ws.Range("D33")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "5.786"
$dCell.Style = $dStyle
$ws.Range("E33").Value = "  -1.83%  "

$dCell = $ws.Range("D34")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "3.609"
$dCell.Style = $dStyle
$ws.Range("E34").Value = "  -0.78%  "

$dCell = $ws.Range("D35")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.02446"
$dCell.Style = $dStyle
$ws.Range("E35").Value = "  +0.10%  "

$dCell = $ws.Range("D36")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.06538"
$dCell.Style = $dStyle
$ws.Range("E36").Value = "  +0.41%  "

$dCell = $ws.Range("D37")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.2158"
$dCell.Style = $dStyle
$ws.Range("E37").Value = "  -1.33%  "

$dCell = $ws.Range("D38")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "8.865"
$dCell.Style = $dStyle
$ws.Range("E38").Value = "  -3.87%  "

$dCell = $ws.Range("D39")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "5.043"
$dCell.Style = $dStyle
$ws.Range("E39").Value = "  +1.25%  "

$dCell = $ws.Range("D40")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.251"
$dCell.Style = $dStyle
$ws.Range("E40").Value = "  -0.64%  "

$dCell = $ws.Range("D41")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.186"
$dCell.Style = $dStyle
$ws.Range("E41").Value = "  -0.17%  "

$dCell = $ws.Range("D42")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.6414"
$dCell.Style = $dStyle
$ws.Range("E42").Value = "  -0.98%  "

$dCell = $ws.Range("D43")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "11.09"
$dCell.Style = $dStyle
$ws.Range("E43").Value = "  -0.60%  "

$dCell = $ws.Range("D44")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.007"
$dCell.Style = $dStyle
$ws.Range("E44").Value = "  -0.63%  "

$dCell = $ws.Range("D45")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "0.6024"
$dCell.Style = $dStyle
$ws.Range("E45").Value = "  -0.56%  "

$dCell = $ws.Range("D46")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "12.96"
$dCell.Style = $dStyle
$ws.Range("E46").Value = "  -0.91%  "

$dCell = $ws.Range("D47")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "3.686"
$dCell.Style = $dStyle
$ws.Range("E47").Value = "  -0.44%  "

$dCell = $ws.Range("D48")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "2.007"
$dCell.Style = $dStyle
$ws.Range("E48").Value = "  -1.11%  "

$dCell = $ws.Range("D49")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.214"
$dCell.Style = $dStyle
$ws.Range("E49").Value = "  -0.35%  "

$dCell = $ws.Range("D50")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "121.94"
$dCell.Style = $dStyle
$ws.Range("E50").Value = "  -0.19%  "

$dCell = $ws.Range("D51")
$dStyle = $dCell.Style
$dCell.NumberFormat = "@"
$dCell.Value = "1.178"
$dCell.Style = $dStyle
$ws.Range("E51").Value = "  -8.92%  "
